$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above row 874, shifting existing rows 874:915 down to 875:916.
$ws.Rows.Item(874).Insert()

# Populate the newly inserted row 874 with the new daily entry.
# Column A holds a date formatted as plain text elsewhere in the sheet, so
# force the Text number format first to keep "2026/02/24" a literal string
# rather than letting Excel auto-convert it to a date serial value.
$ws.Cells.Item(874, 1).NumberFormat = "@"
$ws.Cells.Item(874, 1).Value = "2026/02/24"
$ws.Cells.Item(874, 2).Value = "火"
$ws.Cells.Item(874, 3).Value = 15
$ws.Cells.Item(874, 4).Value = 201
